$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 395 (shifts existing rows 395-425 down to 396-426)
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row 395 with the new record
$ws.Cells.Item(395, 1).Value = 3
$ws.Cells.Item(395, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(395, 3).Value = "Coquimbo"
$ws.Cells.Item(395, 4).Value = 45013
$ws.Cells.Item(395, 5).Value = 5
$ws.Cells.Item(395, 6).Value = 100112001
$ws.Cells.Item(395, 7).Value = "Berenjena"
$ws.Cells.Item(395, 8).Value = "Sin especificar"
$ws.Cells.Item(395, 9).Value = "Primera"
$ws.Cells.Item(395, 10).Value = 130
$ws.Cells.Item(395, 11).Value = 7500
$ws.Cells.Item(395, 12).Value = 8000
$ws.Cells.Item(395, 13).Value = 7769
$ws.Cells.Item(395, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(395, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(395, 16).Value = 129
$ws.Cells.Item(395, 17).Value = 60
$ws.Cells.Item(395, 18).Value = "Hortaliza"

# Match the date cell format used by the rest of column D
$ws.Cells.Item(395, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
